# Auto-generated Excel COM-interop edit script
# Travis County 2018 bg_config.xlsx: interpolation-method / bug-fix recompute.
# Re-orders the factor-variable lists and rewrites the recomputed numeric
# results across all five sheets to match the new run.

$wb = $excel.ActiveWorkbook

# --- "Significant Components": reorder factor-variable lists in column C ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value = "['QHISPC' 'QESL' 'QEDLESHI' 'QNOHLTH' 'PPUNIT' 'QSERV' 'QEXTRCT' 'QFHH'`n 'PERCAP']"
$ws1.Range("C3").Value = "['PERCAP' 'MDHSEVAL' 'QRICH']"
$ws1.Range("C4").Value = "['QRENTER' 'QPOVTY' 'QNOAUTO']"
$ws1.Range("C5").Value = "['QAGEDEP' 'MEDAGE' 'QSSBEN']"

# --- "Loading Factors": re-ordered variable rows with recomputed loadings ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$loadingRows = @(
    @("QHISPC", 0.8230600182487867, 0.3465122430661959, 0.1359332265936176, -0.09949764825394727, -0.1288344852807652),
    @("QESL", 0.7994959189325577, 0.137379140989009, 0.2117197757451876, -0.03267939553219928, -0.2000727788877651),
    @("QEDLESHI", 0.8778539375938658, 0.1961858989009719, 0.1910468548543998, -0.01699417387933021, -0.1019276636881132),
    @("QNOHLTH", 0.6741858254596813, 0.395904438500259, 0.3242825292551726, -0.08306124094097993, -0.1222630312566305),
    @("PPUNIT", 0.7145302567809892, 0.04770843348385883, -0.3603295216403187, -0.0969126445952975, 0.07246473371300596),
    @("QSERV", 0.5274953381277182, 0.3608946343328046, 0.3130727385070557, -0.1262082776137383, -0.1039248313749682),
    @("QEXTRCT", 0.7861441582802978, 0.1131783498334926, 0.08480976403042277, -0.02097332714221105, -0.2028241022156578),
    @("QFHH", 0.5602157009799712, 0.2825588700502056, 0.0438173407653758, -0.06322069614013041, 0.2438533285165921),
    @("PERCAP", 0.474525870673892, 0.7338494729272212, 0.2474036435249251, -0.2041762683047169, 0.03151731764530275),
    @("MDHSEVAL", 0.3004626973795652, 0.7871221079908358, -0.0664155145948064, -0.04968038971114691, 0.02437791650178586),
    @("QRICH", 0.2167683043423544, 0.8465873750415711, 0.3549485669812573, -0.1549872517487549, -0.04401484983639375),
    @("QRENTER", -0.005926325455779259, 0.1780075822013159, 0.770379834680963, -0.4138178503501595, -0.113285003589297),
    @("QPOVTY", 0.3023780538831269, 0.1695709332984764, 0.5567092082697344, -0.2796814542173661, 0.0947606391008271),
    @("QNOAUTO", 0.1139529744537159, 0.06106384662227803, 0.6847406492299022, -0.05243788690298962, 0.007215022956936492),
    @("QAGEDEP", -0.01897725427186348, -0.1465327647892834, -0.1088299637021307, 0.6889723577247717, 0.6048116249732312),
    @("MEDAGE", -0.3203145700577555, -0.2374924602301702, -0.358121827837174, 0.7551459454149352, -0.03604536918156718),
    @("QSSBEN", 0.006717324505023358, -0.04185582980288473, -0.1566751782671676, 0.8301134344857622, 0.08972562297410745),
    @("QFEMLBR", -0.2650015955920006, 0.08497030028627561, 0.005596367281132848, -0.04604931303354207, 0.7314609899884509),
    @("QFEMALE", -0.01889530618167813, -0.05607244865644413, 0.002629353904375671, 0.1873923070417907, 0.8727849552670841)
)
for ($i = 0; $i -lt $loadingRows.Count; $i++) {
    $r = 2 + $i
    $row = $loadingRows[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
}

# --- "All Refactor Variances": recomputed SS Loadings / variances (F2 block, cols N:R) ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$refactorRows = @(
    @(4.944807421852724, 2.838952307033347, 2.400259439223039, 2.313501490430895, 1.905085930346282, 1.541385597189871, 1.246987809407878, 4.751932245886838, 3.118803346904831, 2.379793098020169, 2.141113004780542, 1.886066872631294, 4.902928251948794, 2.588309621185193, 2.170494330712241, 2.139559861103693, 1.895944383586705),
    @(0.183141015624175, 0.1051463817419758, 0.08889849774900146, 0.08568524038632945, 0.0705587381609734, 0.05708835545147672, 0.04618473368177325, 0.2262824878993732, 0.1485144450907062, 0.1133234808581033, 0.1019577621324068, 0.08981270822053783, 0.258048855365726, 0.1362268221676418, 0.1142365437216969, 0.1126084137422996, 0.09978654650456342),
    @(0.183141015624175, 0.2882873973661508, 0.3771858951151522, 0.4628711355014817, 0.5334298736624551, 0.5905182291139318, 0.636702962795705, 0.2262824878993732, 0.3747969329900795, 0.4881204138481828, 0.5900781759805895, 0.6798908842011273, 0.258048855365726, 0.3942756775333678, 0.5085122212550647, 0.6211206349973643, 0.7209071815019277),
    @(0.2876396472540655, 0.1651419702529537, 0.1396231884310012, 0.1345764750490454, 0.1108189254391975, 0.08966246238403999, 0.0725373311896968, 0.332821770606993, 0.2184386473503185, 0.1666789237676845, 0.1499619490445254, 0.1320987092304788, 0.3579501799775537, 0.1889658275893836, 0.1584622079692675, 0.1562037619152203, 0.1384180225485749)
)
for ($i = 0; $i -lt $refactorRows.Count; $i++) {
    $r = 2 + $i
    $row = $refactorRows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws3.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}

# --- "Final Variances": recomputed SS Loadings / variances ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$finalRows = @(
    @(4.902928251948794, 2.588309621185193, 2.170494330712241, 2.139559861103693, 1.895944383586705),
    @(0.258048855365726, 0.1362268221676418, 0.1142365437216969, 0.1126084137422996, 0.09978654650456342),
    @(0.258048855365726, 0.3942756775333678, 0.5085122212550647, 0.6211206349973643, 0.7209071815019277),
    @(0.3579501799775537, 0.1889658275893836, 0.1584622079692675, 0.1562037619152203, 0.1384180225485749)
)
for ($i = 0; $i -lt $finalRows.Count; $i++) {
    $r = 2 + $i
    $row = $finalRows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws4.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
}

# --- "Included and Excluded": reordered include list ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value = "[['QHISPC', 'QESL', 'QEDLESHI', 'QNOHLTH', 'PPUNIT', 'QSERV', 'QEXTRCT', 'QFHH', 'PERCAP', 'MDHSEVAL', 'QRICH', 'QRENTER', 'QPOVTY', 'QNOAUTO', 'QAGEDEP', 'MEDAGE', 'QSSBEN', 'QFEMLBR', 'QFEMALE']]"

